$d = $word.ActiveDocument

$d.Content.Find.Execute(" SE*2.33 = ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " SE*2.33 = ", 2)

$d.Content.Find.Execute("APPROXIMATELY 98% ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "APPROXIMATELY 98% ", 2)

$d.Content.Find.Execute("of sample means fall w/in some interval", $true, $false, $false, $false, $false,
                         $true, 1, $false, "of sample means fall w/in some interval", 2)
